# "Generate Report for Handback"
#
# - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   on every sheet that carries a Status column (Overview + the two language
#   sheets).
# - Each language sheet (zh-cn / de-de) grows two new report columns:
#     F = Latest Target File   (the translated file, same name as the source .md)
#     G = Latest Handback File (the handback .xlf)
#   for both data rows, each rendered as a hyperlink exactly like the
#   existing Latest Handoff File / Source File Name columns.
# - de-de's "Latest Handback DateTime" (column H) is stamped with the
#   actual handback time; zh-cn's handback time is refreshed too.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Status column: "Ready for handoff" -> "Handed back: in sync with en-US" ---
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# --- zh-cn: add Latest Target File (F) / Latest Handback File (G) hyperlinks ---
$zhXlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/67091fba48c9a781dc0ad1cc2aa12003b4c162a7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/" + $zhXlf

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/42f3aa706d26137687e17f337678e094984159b1/e2e/a.md", "", "", "a.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G2"), $zhXlfUrl, "", "", $zhXlf)
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/42f3aa706d26137687e17f337678e094984159b1/e2e/a.md", "", "", "a.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("G3"), $zhXlfUrl, "", "", $zhXlf)

# Refresh the zh-cn handback timestamp (Latest Handback DateTime, column H)
$wsZhCn.Range("H2").Value = "2016-03-20 16:30:51"
$wsZhCn.Range("H3").Value = "2016-03-20 16:30:51"

# --- de-de: add Latest Target File (F) / Latest Handback File (G) hyperlinks ---
$deXlf = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a178c833e5ad4485374a4ab13493d20481a5df25/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/" + $deXlf

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/42f3aa706d26137687e17f337678e094984159b1/e2e/a.md", "", "", "a.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G2"), $deXlfUrl, "", "", $deXlf)
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/42f3aa706d26137687e17f337678e094984159b1/e2e/a.md", "", "", "a.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("G3"), $deXlfUrl, "", "", $deXlf)

# Refresh the de-de handback timestamp (Latest Handback DateTime, column H)
$wsDeDe.Range("H2").Value = "2016-03-20 16:30:57"
$wsDeDe.Range("H3").Value = "2016-03-20 16:30:57"

Write-Output "Handback report generated."
